$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) "Variable Instalments"
# column before the existing "Late" column (N), pushing Late/Outstanding
# (heading)/Outstanding one column to the right. ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Column N currently holds "Late"; insert a blank column in front of it so
# the existing Late/heading/Outstanding columns shift right by one.
$wsSchedule.Columns("N:N").Insert()

# The newly inserted column should carry the same width as its neighbour
# (column M, "In Advance") rather than the sheet default.
$wsSchedule.Columns("N:N").ColumnWidth = $wsSchedule.Columns("M:M").ColumnWidth

# Activate the "Repayment schedule" sheet and move the selection to S7 -
# this makes it the workbook's active tab (and clears "tabSelected" /
# the previous scrolled position on whichever sheet had it before, i.e.
# "Transactions").
$wsSchedule.Activate()
$wsSchedule.Range("S7").Select()
